$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update shared string description for observed variable (B1 header)
$ws.Range("B1").Value = "nom_w_pc_obs"

# Update raw variable values in column B (rows 2-104), including new row 104
$ws.Range("B2").Value = -0.001963394164474025
$ws.Range("B3").Value = -0.002327504873449959
$ws.Range("B4").Value = -0.008373402913548197
$ws.Range("B5").Value = -0.0002334146387425184
$ws.Range("B6").Value = -0.004974087501678437
$ws.Range("B7").Value = -0.002873876727537014
$ws.Range("B8").Value = 0.001872250636571593
$ws.Range("B9").Value = -0.00228381552655392
$ws.Range("B10").Value = -0.002278336726335276
$ws.Range("B11").Value = -0.0004924139046595666
$ws.Range("B12").Value = 0.006035274326563464
$ws.Range("B13").Value = 0.007769347920333403
$ws.Range("B14").Value = 0.002415724089875573
$ws.Range("B15").Value = 0.00118931654183134
$ws.Range("B16").Value = 0.001565370138131392
$ws.Range("B17").Value = 0.001790707291433644
$ws.Range("B18").Value = -0.00694435872936372
$ws.Range("B19").Value = -0.002425624942349477
$ws.Range("B20").Value = -0.001008517241987528
$ws.Range("B21").Value = 0.001184001079845579
$ws.Range("B22").Value = 0.006224470727149073
$ws.Range("B23").Value = -0.0008061487366662068
$ws.Range("B24").Value = -0.003361268411086568
$ws.Range("B25").Value = -0.004484875849591119
$ws.Range("B26").Value = -0.01132433911920276
$ws.Range("B27").Value = -0.01001255424457459
$ws.Range("B28").Value = -0.001729372575774657
$ws.Range("B29").Value = 0.004175660861490085
$ws.Range("B30").Value = -0.001575335022303712
$ws.Range("B31").Value = 0.007002297589531081
$ws.Range("B32").Value = 0.003444667621483974
$ws.Range("B33").Value = -0.002518721722083173
$ws.Range("B34").Value = -0.004773396740810676
$ws.Range("B35").Value = 0.009413578462240307
$ws.Range("B36").Value = 0.0004995165982373123
$ws.Range("B37").Value = 0.001965463999947098
$ws.Range("B38").Value = 0.01060139104450275
$ws.Range("B39").Value = -0.0008800727788100815
$ws.Range("B40").Value = 0.007358359036720559
$ws.Range("B41").Value = -0.005788029093942874
$ws.Range("B42").Value = -0.007631792155466721
$ws.Range("B43").Value = -0.004671021829008942
$ws.Range("B44").Value = -0.003510703031341253
$ws.Range("B45").Value = -0.004611063363654272
$ws.Range("B46").Value = 0.002955948045206061
$ws.Range("B47").Value = -0.009215908231850908
$ws.Range("B48").Value = -0.009272178687502755
$ws.Range("B49").Value = -0.003046135161522259
$ws.Range("B50").Value = -0.0003963008488961223
$ws.Range("B51").Value = -0.001672419281086079
$ws.Range("B52").Value = -0.001460118188855108
$ws.Range("B53").Value = 0.000322224279866376
$ws.Range("B54").Value = 0.001359508163468234
$ws.Range("B55").Value = 0.002622025574123513
$ws.Range("B56").Value = -0.001234995625449209
$ws.Range("B57").Value = -0.004240339388076789
$ws.Range("B58").Value = -0.004141632451587704
$ws.Range("B59").Value = 0.000610156620693289
$ws.Range("B60").Value = 0.004047255672081682
$ws.Range("B61").Value = 0.01040248382502396
$ws.Range("B62").Value = 0.0145316991866177
$ws.Range("B63").Value = 0.007954394128898014
$ws.Range("B64").Value = 0.008539955473361269
$ws.Range("B65").Value = -0.002543627116747119
$ws.Range("B66").Value = 0.008302946294300263
$ws.Range("B67").Value = -0.005348561326916364
$ws.Range("B68").Value = 0.001276359773586241
$ws.Range("B69").Value = 0.01323850717418746
$ws.Range("B70").Value = 0.01415004597231243
$ws.Range("B71").Value = -0.005583624382106071
$ws.Range("B72").Value = 0.01098340024974537
$ws.Range("B73").Value = -0.002811435893628575
$ws.Range("B74").Value = 0.01256398819791785
$ws.Range("B75").Value = -0.007508911314378432
$ws.Range("B76").Value = -0.003642373631874308
$ws.Range("B77").Value = 0.003188310844850489
$ws.Range("B78").Value = 0.003596726868077233
$ws.Range("B79").Value = 0.001951000246955345
$ws.Range("B80").Value = -0.002709270498215394
$ws.Range("B81").Value = -0.004752789629131809
$ws.Range("B82").Value = 0.001053642528424417
$ws.Range("B83").Value = 0.01035110686425354
$ws.Range("B84").Value = 0.004293346537803751
$ws.Range("B85").Value = -0.00115305671897727
$ws.Range("B86").Value = -0.01432597766720594
$ws.Range("B87").Value = 0.002015120475835603
$ws.Range("B88").Value = 0.004727584164920291
$ws.Range("B89").Value = -0.003331277738278505
$ws.Range("B90").Value = -0.004862589571685204
$ws.Range("B91").Value = -0.003628658572890836
$ws.Range("B92").Value = -0.00001960821037370192
$ws.Range("B93").Value = -0.006343162706074861
$ws.Range("B94").Value = 0.002690668643429795
$ws.Range("B95").Value = -0.008187966377292711
$ws.Range("B96").Value = -0.008106612773870007
$ws.Range("B97").Value = 0.01476256641352053
$ws.Range("B98").Value = -0.006982890863131128
$ws.Range("B99").Value = -0.006430824044111862
$ws.Range("B100").Value = -0.00009805544284256784
$ws.Range("B101").Value = 0.00698338133195929
$ws.Range("B102").Value = 0.00792942162362931
$ws.Range("B103").Value = -0.01212899702642564
$ws.Range("B104").Value = -0.002873431204955041
